$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D stays text (values look numeric but must remain strings,
# matching the "Price" column formatting used throughout the sheet).
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '39.466.82'
$ws.Range("E2").Value = '  +1.83%  '

# Row 3
$ws.Range("D3").Value = '2.172.79'
$ws.Range("E3").Value = '  +3.90%  '

# Row 4
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
$ws.Range("D5").Value = '229.43'
$ws.Range("E5").Value = '  +0.49%  '

# Row 6
$ws.Range("D6").Value = '0.622'
$ws.Range("E6").Value = '  +1.16%  '

# Row 7
$ws.Range("D7").Value = '62.91'
$ws.Range("E7").Value = '  +3.33%  '

# Row 8
$ws.Range("D8").Value = '0.997'
$ws.Range("E8").Value = '  -0.26%  '

# Row 9
$ws.Range("D9").Value = '0.393'
$ws.Range("E9").Value = '  +2.13%  '

# Row 10
$ws.Range("D10").Value = '0.0854'
$ws.Range("E10").Value = '  +1.89%  '

# Row 11
$ws.Range("E11").Value = '  -0.39%  '

# Row 12
$ws.Range("D12").Value = '16.09'
$ws.Range("E12").Value = '  +7.43%  '

# Row 13
$ws.Range("D13").Value = '2.493.24'
$ws.Range("E13").Value = '  +3.94%  '

# Row 14
$ws.Range("D14").Value = '22.33'
$ws.Range("E14").Value = '  +1.91%  '

# Row 15
$ws.Range("D15").Value = '0.821'
$ws.Range("E15").Value = '  +3.13%  '

# Row 16
$ws.Range("D16").Value = '5.54'
$ws.Range("E16").Value = '  +0.97%  '

# Row 17
$ws.Range("D17").Value = '2.160.71'
$ws.Range("E17").Value = '  +3.54%  '

# Row 18
$ws.Range("D18").Value = '39.346.09'
$ws.Range("E18").Value = '  +1.66%  '

# Row 19
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").Value = '6.16'
$ws.Range("E19").Value = '  +2.08%  '

# Row 20
$ws.Range("B20").Value = 'Litecoin'
$ws.Range("C20").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D20").Value = '72.13'
$ws.Range("E20").Value = '  +0.60%  '

# Row 21
$ws.Range("D21").Value = '0.0₃0852'
$ws.Range("E21").Value = '  +1.85%  '

# Row 22
$ws.Range("D22").Value = '228.15'
$ws.Range("E22").Value = '  +0.76%  '

# Row 23
$ws.Range("E23").Value = '  +0.16%  '

# Row 24
$ws.Range("E24").Value = '  +1.22%  '

# Row 25
$ws.Range("D25").Value = '2.34'
$ws.Range("E25").Value = '  +0.04%  '

# Row 26
$ws.Range("D26").Value = '9.72'
$ws.Range("E26").Value = '  +2.95%  '

# Row 27
$ws.Range("D27").Value = '170.94'
$ws.Range("E27").Value = '  +0.13%  '

# Row 28
$ws.Range("D28").Value = '0.137'
$ws.Range("E28").Value = '  -0.23%  '

# Row 29
$ws.Range("D29").Value = '19.71'
$ws.Range("E29").Value = '  +2.84%  '

# Row 30
$ws.Range("D30").Value = '1.40'
$ws.Range("E30").Value = '  -2.00%  '

# Row 31
$ws.Range("E31").Value = '  +10.01%  '

# Row 32
$ws.Range("E32").Value = '  +0.66%  '

# Row 33
$ws.Range("D33").Value = '4.62'
$ws.Range("E33").Value = '  +2.64%  '

# Row 34
$ws.Range("D34").Value = '4.83'
$ws.Range("E34").Value = '  +2.58%  '

# Row 35
$ws.Range("D35").Value = '7.22'
$ws.Range("E35").Value = '  +12.76%  '

# Row 36
$ws.Range("D36").Value = '0.0620'
$ws.Range("E36").Value = '  +0.83%  '

# Row 37
$ws.Range("E37").Value = '  +1.85%  '

# Row 38
$ws.Range("D38").Value = '3.57'
$ws.Range("E38").Value = '  +1.23%  '

# Row 39
$ws.Range("D39").Value = '0.998'
$ws.Range("E39").Value = '  -0.14%  '

# Row 40
$ws.Range("D40").Value = '18.23'
$ws.Range("E40").Value = '  -0.27%  '

# Row 41
$ws.Range("D41").Value = '0.0231'
$ws.Range("E41").Value = '  +3.86%  '

# Row 42
$ws.Range("D42").Value = '103.47'
$ws.Range("E42").Value = '  +2.48%  '

# Row 43
$ws.Range("D43").Value = '1.536.24'
$ws.Range("E43").Value = '  -0.19%  '

# Row 44
$ws.Range("D44").Value = '1.20'
$ws.Range("E44").Value = '  +6.89%  '

# Row 45
$ws.Range("D45").Value = '1.11'
$ws.Range("E45").Value = '  +7.30%  '

# Row 46
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").Value = '7.83'
$ws.Range("E46").Value = '  +1.65%  '

# Row 47
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = '0.0923'
$ws.Range("E47").Value = '  -0.25%  '

# Row 48
$ws.Range("B48").Value = 'HuobiToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D48").Value = '2.80'
$ws.Range("E48").Value = '  -0.51%  '

# Row 49
$ws.Range("E49").Value = '  +1.40%  '

# Row 50
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.376.99'
$ws.Range("E50").Value = '  +3.95%  '

# Row 51
$ws.Range("B51").Value = 'MXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D51").Value = '2.99'
$ws.Range("E51").Value = '  +0.56%  '
